# edit.ps1 — applies the "201 Created" -> "200 Created" fix to the
# "Test Cases" sheet's Expected_Result (J) column, and clears the
# leftover empty Notes (F) cells on the "Execution Tracking" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. "Test Cases" sheet: Expected_Result column (J) — 201 -> 200
# ---------------------------------------------------------------
$wsTestCases = $wb.Worksheets.Item("Test Cases")

$rowsToFix = @(2, 3, 4, 5, 14, 15, 16, 17)
foreach ($r in $rowsToFix) {
    $cell = $wsTestCases.Cells.Item($r, 10)  # column J = 10
    $current = $cell.Formula
    $cell.Value = $current -replace '^201 Created', '200 Created'
}

# ---------------------------------------------------------------
# 2. "Execution Tracking" sheet: drop the empty inline-string cells
#    left behind in column F (Notes) for rows 2-32.
# ---------------------------------------------------------------
$wsExecTracking = $wb.Worksheets.Item("Execution Tracking")

for ($r = 2; $r -le 32; $r++) {
    $wsExecTracking.Cells.Item($r, 6).ClearContents()  # column F = 6
}
